$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CaseTypeTab")

# Insert a new column before the existing "TabFieldDisplayOrder" column (I)
# so the new UserRole column becomes column I, shifting the rest right.
$ws.Columns.Item(9).Insert()

# New header / hint cell for the UserRole column (row 2 = descriptions row)
$ws.Cells.Item(2, 9).Value = "MaxLength: 100. No entry for role means no role restriction for that tab. Enter role on a single row per tab"

# New column title (row 3 = column names row)
$ws.Cells.Item(3, 9).Value = "UserRole"

# Make CaseTypeTab the active/selected sheet (moves tabSelected from FixedLists)
$ws.Activate()
$ws.Range("I3").Select()
